# Applies "semana 42 de 2024" updates to the poisson.xlsx sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 0

# Row 5
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 0.13

# Row 6
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 58

# Row 7
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0.06

# Row 11
$ws.Range("C11").Value = 46
$ws.Range("D11").Value = 25

# Row 13
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 0

# Row 14
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 7
$ws.Range("E14").Value = 0.06

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0.37

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0.02

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0.37

# Row 21
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 0

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("E22").Value = 0.37

# Row 24
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 0.27

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0.03

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0.37

# Row 29
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0

# Row 31
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 0.27

# Row 34
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 0.27

# Row 35
$ws.Range("C35").Value = 5
$ws.Range("D35").Value = 4
$ws.Range("E35").Value = 0.18

# Row 36
$ws.Range("C36").Value = 8
$ws.Range("D36").Value = 2
$ws.Range("E36").Value = 0.01

# Row 37
$ws.Range("C37").Value = 9
$ws.Range("D37").Value = 9
$ws.Range("E37").Value = 0.13
